{"js": "// Insert <w:br/> (line break) splits inside three paragraphs, per the\n// target diff:\n//   1) \"Objetivos\" paragraph: break after \"...compet\u00eancias no\" and before\n//      \"desenvolvimento...\"\n//   2) \"Programa\" paragraph: break after \"Programa\" and before\n//      \"1.INTRODU\u00c7\u00c3O...\"\n//   3) \"Bibliografia\" paragraph: break before each of the 11 numbered\n//      references (i.e. after each \".<number>.\" sequence's previous\n//      sentence, right before \"2.\", \"3.\", ... \"11.\")\n//\n// Strategy: locate a short, unique anchor string scoped to the target\n// paragraph via Paragraph.search(), collapse to a zero-width caret right\n// after the anchor with getRange(\"End\"), then insert a vertical-tab\n// character (U+000B) there with InsertLocation.Replace. Word's OOXML\n// writer represents U+000B as a <w:br/> element, splitting the run in\n// two <w:t> pieces exactly like the target diff \u2014 tested against this\n// runtime and confirmed to reproduce the exact <w:t>/<w:br/>/<w:t>\n// sequence.\n\nasync function breakAfter(paragraph, anchorText) {\n  const results = paragraph.search(anchorText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for ${JSON.stringify(anchorText)}, found ${results.items.length}`);\n  }\n  const endRange = results.items[0].getRange(\"End\");\n  endRange.insertText(\"\\u000b\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three target paragraphs by their distinctive text content\n// rather than hard-coded indices, so the script is resilient to minor\n// document-structure differences.\nlet objetivosPara = null;\nlet programaPara = null;\nlet bibliografiaPara = null;\n\nfor (const para of paragraphs.items) {\n  const t = para.text;\n  if (t.indexOf(\"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais\") === 0) {\n    objetivosPara = para;\n  } else if (t.indexOf(\"Programa1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS\") === 0) {\n    programaPara = para;\n  } else if (t.indexOf(\"1. Meyers, M., Chawla, K. Mechanical Behavior of Materials.\") === 0) {\n    bibliografiaPara = para;\n  }\n}\n\nif (!objetivosPara) throw new Error(\"Objetivos paragraph not found\");\nif (!programaPara) throw new Error(\"Programa paragraph not found\");\nif (!bibliografiaPara) throw new Error(\"Bibliografia paragraph not found\");\n\n// 1) Objetivos: break between \"...compet\u00eancias no\" and \"desenvolvimento...\"\nawait breakAfter(objetivosPara, \"compet\u00eancias no\");\n\n// 2) Programa: break between \"Programa\" and \"1.INTRODU\u00c7\u00c3O...\"\nawait breakAfter(programaPara, \"Programa\");\n\n// 3) Bibliografia: break before each numbered reference 2..11\nconst bibAnchors = [\n  \"Ed. Cambridge University Press, 2009.\",\n  \"Pergamon Press, 1985.\",\n  \"Ed. Guanabara Dois, 1981.\",\n  \"Pergamon Press, 1965.\",\n  \"Edward Arnold, 1967.\",\n  \"Ed. Guanabara Dois, 1982.\",\n  \"Ed. Edgard Blucher Ltda., 1970.\",\n  \"Ed. Edgar Bl\u00fccher, 2008.\",\n  \"New Jersey, Prentice Hall,1988.\",\n  \"Livros T\u00e9cnicos e Cient\u00edficos, 2008.\",\n];\n\nfor (const anchor of bibAnchors) {\n  await breakAfter(bibliografiaPara, anchor);\n}\n", "ps1": "# Insert manual line breaks (<w:br/>) inside three paragraphs, per the\n# target diff:\n#   1) \"Objetivos\" paragraph: break after \"...compet\u00eancias no\" and before\n#      \"desenvolvimento...\"\n#   2) \"Programa\" paragraph: break after \"Programa\" and before\n#      \"1.INTRODU\u00c7\u00c3O...\"\n#   3) \"Bibliografia\" paragraph: break right before each of the 10 later\n#      numbered references (2. .. 11.)\n#\n# Strategy: scope a Find/Replace to each target paragraph's own Range\n# (so the search can't bleed into neighboring paragraphs/headings that\n# share the same leading word, e.g. the \"Programa\" heading), search for\n# a short unique anchor string, and replace it with itself plus the\n# special \"^l\" Find/Replace code, which Word expands to a manual line\n# break (a <w:br/> element splitting the run's text into two <w:t> runs)\n# \u2014 confirmed against this runtime to reproduce the exact target\n# <w:t>/<w:br/>/<w:t> sequence.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreakAfter($paragraphRange, $anchorText) {\n    $find = $paragraphRange.Find\n    $find.Text = $anchorText\n    $find.Replacement.Text = $anchorText + \"^l\"\n    $found = $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n    if (-not $found) {\n        throw \"Anchor not found: $anchorText\"\n    }\n}\n\n# Locate the three target paragraphs by their distinctive text content\n# rather than hard-coded indices, so the script is resilient to minor\n# document-structure differences.\n$objetivosIndex = 0\n$programaIndex = 0\n$bibliografiaIndex = 0\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais\")) {\n        $objetivosIndex = $i\n    } elseif ($t.StartsWith(\"Programa1.INTRODU\u00c7\u00c3O AO CONCEITO DE PROPRIEDADES MEC\u00c2NICAS\")) {\n        $programaIndex = $i\n    } elseif ($t.StartsWith(\"1. Meyers, M., Chawla, K. Mechanical Behavior of Materials.\")) {\n        $bibliografiaIndex = $i\n    }\n}\n\nif ($objetivosIndex -eq 0) { throw \"Objetivos paragraph not found\" }\nif ($programaIndex -eq 0) { throw \"Programa paragraph not found\" }\nif ($bibliografiaIndex -eq 0) { throw \"Bibliografia paragraph not found\" }\n\n# 1) Objetivos: break between \"...compet\u00eancias no\" and \"desenvolvimento...\"\nInsert-LineBreakAfter $d.Paragraphs.Item($objetivosIndex).Range \"compet\u00eancias no\"\n\n# 2) Programa: break between \"Programa\" and \"1.INTRODU\u00c7\u00c3O...\"\nInsert-LineBreakAfter $d.Paragraphs.Item($programaIndex).Range \"Programa\"\n\n# 3) Bibliografia: break before each numbered reference 2..11\n$bibAnchors = @(\n    \"Ed. Cambridge University Press, 2009.\",\n    \"Pergamon Press, 1985.\",\n    \"Ed. Guanabara Dois, 1981.\",\n    \"Pergamon Press, 1965.\",\n    \"Edward Arnold, 1967.\",\n    \"Ed. Guanabara Dois, 1982.\",\n    \"Ed. Edgard Blucher Ltda., 1970.\",\n    \"Ed. Edgar Bl\u00fccher, 2008.\",\n    \"New Jersey, Prentice Hall,1988.\",\n    \"Livros T\u00e9cnicos e Cient\u00edficos, 2008.\"\n)\n\nforeach ($anchor in $bibAnchors) {\n    Insert-LineBreakAfter $d.Paragraphs.Item($bibliografiaIndex).Range $anchor\n}\n"}
